$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3308
$ws.Range("F5").Value = 1356
$ws.Range("F7").Value = 3823
$ws.Range("F9").Value = 191
$ws.Range("F11").Value = 8450
$ws.Range("F12").Value = 8450
$ws.Range("F13").Value = 458
$ws.Range("F15").Value = 130
$ws.Range("F22").Value = 10632
$ws.Range("F23").Value = 10632
$ws.Range("F27").Value = 140
$ws.Range("F38").Value = 2111
$ws.Range("F40").Value = 4064
$ws.Range("F41").Value = 44
$ws.Range("F48").Value = 330
$ws.Range("F49").Value = 289
$ws.Range("F51").Value = 107

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 3
$ws.Range("F10").Value = 30
$ws.Range("F20").Value = 27
$ws.Range("F22").Value = 41

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 3308
$ws.Range("F8").Value = 1356
$ws.Range("F10").Value = 3823
$ws.Range("F15").Value = 191
$ws.Range("F16").Value = 8450
$ws.Range("F17").Value = 458
$ws.Range("F19").Value = 130
$ws.Range("F25").Value = 10632
$ws.Range("F29").Value = 140
$ws.Range("F41").Value = 2111
$ws.Range("F44").Value = 45
$ws.Range("F48").Value = 330
$ws.Range("F49").Value = 289
$ws.Range("F51").Value = 107
